$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 6 de Mayo de 2020 a las 11:03"

# Swap country-name labels for the three rank-changed row pairs
$ws.Range("A42").Value = "Filipinas"
$ws.Range("A43").Value = "Dinamarca"
$ws.Range("A67").Value = "Oman"
$ws.Range("A68").Value = "Armenia"
$ws.Range("A205").Value = "Montserrat"
$ws.Range("A206").Value = "Seychelles"

# Update numeric statistics cells per the source refresh
$ws.Range("B32").Value = 16314
$ws.Range("C32").Value = 25
$ws.Range("D32").Value = 10527
$ws.Range("E32").Value = 5549
$ws.Range("F32").Value = 90
$ws.Range("B33").Value = 15684
$ws.Range("C33").Value = 34
$ws.Range("D33").Value = 13639
$ws.Range("E33").Value = 1437
$ws.Range("F33").Value = 97
$ws.Range("G33").Value = 2
$ws.Range("H33").Value = 608
$ws.Range("B36").Value = 14647
$ws.Range("C36").Value = 216
$ws.Range("E36").Value = 9269
$ws.Range("G36").Value = 7
$ws.Range("H36").Value = 723
$ws.Range("B42").Value = 10004
$ws.Range("C42").Value = 320
$ws.Range("D42").Value = 1506
$ws.Range("E42").Value = 7840
$ws.Range("F42").Value = 31
$ws.Range("G42").Value = 21
$ws.Range("H42").Value = 658
$ws.Range("B43").Value = 9938
$ws.Range("C43").Value = 117
$ws.Range("D43").Value = 7296
$ws.Range("E43").Value = 2139
$ws.Range("F43").Value = 49
$ws.Range("H43").Value = 503
$ws.Range("B53").Value = 6428
$ws.Range("C53").Value = 45
$ws.Range("D53").Value = 4702
$ws.Range("E53").Value = 1619
$ws.Range("F53").Value = 22
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = 107
$ws.Range("B67").Value = 2903
$ws.Range("C67").Value = 168
$ws.Range("D67").Value = 888
$ws.Range("E67").Value = 2002
$ws.Range("F67").Value = 17
$ws.Range("H67").Value = 13
$ws.Range("B68").Value = 2782
$ws.Range("C68").Value = 163
$ws.Range("D68").Value = 1135
$ws.Range("E68").Value = 1607
$ws.Range("F68").Value = 10
$ws.Range("H68").Value = 40
$ws.Range("D72").Value = 1547
$ws.Range("E72").Value = 660
$ws.Range("B81").Value = 1713
$ws.Range("C81").Value = 2
$ws.Range("D81").Value = 264
$ws.Range("E81").Value = 1394
$ws.Range("F81").Value = 4
$ws.Range("B88").Value = 1428
$ws.Range("C88").Value = 5
$ws.Range("D88").Value = 718
$ws.Range("E88").Value = 662
$ws.Range("G88").Value = 2
$ws.Range("H88").Value = 48
$ws.Range("D92").Value = 932
$ws.Range("E92").Value = 105
$ws.Range("F92").Value = 2
$ws.Range("D101").Value = 215
$ws.Range("E101").Value = 547
$ws.Range("D130").Value = 261
$ws.Range("E130").Value = 55
$ws.Range("D205").Value = 7
$ws.Range("F205").Value = 1
$ws.Range("H205").Value = 1
$ws.Range("D206").Value = 8
$ws.Range("F206").Value = 0
$ws.Range("H206").Value = 0
